# Insert a new data row before the existing row 346 (shifting rows 346:435
# down to 347:436), then populate the newly inserted row 346 with a fresh
# "Choclo" price record for "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 346; existing rows 346:435 shift to 347:436
$ws.Rows.Item(346).Insert()

# Populate the new row 346
$ws.Range("A346").Value2 = 10
$ws.Range("B346").Value2 = "Vega Modelo de Temuco"
$ws.Range("C346").Value2 = "La Araucanía"
$ws.Range("D346").Value2 = 44642
$ws.Range("E346").Value2 = 9
$ws.Range("F346").Value2 = 100112024
$ws.Range("G346").Value2 = "Choclo"
$ws.Range("H346").Value2 = "Dulce o Americano"
$ws.Range("I346").Value2 = "Primera"
$ws.Range("J346").Value2 = 15000
$ws.Range("K346").Value2 = 120
$ws.Range("L346").Value2 = 130
$ws.Range("M346").Value2 = 123
$ws.Range("N346").Value2 = "$/unidad"
$ws.Range("O346").Value2 = "Región de La Araucanía"
$ws.Range("P346").Value2 = 123
$ws.Range("Q346").Value2 = 1
$ws.Range("R346").Value2 = "Hortaliza"

# Keep the date-format style used by the rest of column D
$ws.Range("D346").NumberFormat = $ws.Range("D347").NumberFormat
